# Refresh the crypto price/volume snapshot (values scraped from coinranking.com).
# Source data is plain text (prices keep trailing zeros / locale-style dotted
# thousands separators like "28.314.53", volumes are "  +2.24%  " with
# padding spaces) so every write goes through Set-TextCell, which forces the
# cell to Text format first for any value Excel's .Value setter would
# otherwise auto-coerce into a Number (e.g. "22.00" -> 22, "1.140" -> 1.14).
# Values that are already unambiguous as text (two-dot strings, %-strings,
# coin names/URLs) skip the NumberFormat step since plain assignment already
# keeps them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Row, $Col, [string]$Text, [bool]$ForceText = $false)
    $cell = $Sheet.Cells.Item($Row, $Col)
    if ($ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $Text
}

# Row 2
Set-TextCell $ws 2 4 '28.314.53' $false
Set-TextCell $ws 2 5 '  +2.24%  ' $false

# Row 3
Set-TextCell $ws 3 4 '1.814.31' $false
Set-TextCell $ws 3 5 '  +3.37%  ' $false

# Row 4
Set-TextCell $ws 4 5 '  -0.05%  ' $false

# Row 5
Set-TextCell $ws 5 4 '325.76' $true
Set-TextCell $ws 5 5 '  +0.40%  ' $false

# Row 6
Set-TextCell $ws 6 4 '0.9991' $true
Set-TextCell $ws 6 5 '  +0.05%  ' $false

# Row 7
Set-TextCell $ws 7 4 '0.4371' $true
Set-TextCell $ws 7 5 '  +2.04%  ' $false

# Row 8
Set-TextCell $ws 8 4 '0.3668' $true
Set-TextCell $ws 8 5 '  +0.76%  ' $false

# Row 9
Set-TextCell $ws 9 4 '44.66' $true
Set-TextCell $ws 9 5 '  -1.35%  ' $false

# Row 10
Set-TextCell $ws 10 4 '0.07676' $true
Set-TextCell $ws 10 5 '  +2.45%  ' $false

# Row 11
Set-TextCell $ws 11 4 '1.140' $true
Set-TextCell $ws 11 5 '  +1.48%  ' $false

# Row 12
Set-TextCell $ws 12 4 '0.9991' $true
Set-TextCell $ws 12 5 '  +0.04%  ' $false

# Row 13
Set-TextCell $ws 13 4 '22.00' $true
Set-TextCell $ws 13 5 '  +1.60%  ' $false

# Row 14
Set-TextCell $ws 14 4 '6.308' $true
Set-TextCell $ws 14 5 '  +2.58%  ' $false

# Row 15
Set-TextCell $ws 15 4 '7.491' $true
Set-TextCell $ws 15 5 '  +3.15%  ' $false

# Row 16
Set-TextCell $ws 16 4 '1.821.79' $false
Set-TextCell $ws 16 5 '  +4.25%  ' $false

# Row 17
Set-TextCell $ws 17 4 '95.36' $true
Set-TextCell $ws 17 5 '  +8.74%  ' $false

# Row 18
Set-TextCell $ws 18 5 '  +0.77%  ' $false

# Row 19
Set-TextCell $ws 19 4 '0.06501' $true
Set-TextCell $ws 19 5 '  +4.81%  ' $false

# Row 20
Set-TextCell $ws 20 4 '0.9987' $true
Set-TextCell $ws 20 5 '  +0.00%  ' $false

# Row 21
Set-TextCell $ws 21 4 '17.39' $true
Set-TextCell $ws 21 5 '  +1.85%  ' $false

# Row 22
Set-TextCell $ws 22 4 '6.251' $true
Set-TextCell $ws 22 5 '  +1.50%  ' $false

# Row 23
Set-TextCell $ws 23 4 '28.306.47' $false
Set-TextCell $ws 23 5 '  +2.20%  ' $false

# Row 24
Set-TextCell $ws 24 4 '11.56' $true
Set-TextCell $ws 24 5 '  -1.36%  ' $false

# Row 25
Set-TextCell $ws 25 4 '2.115' $true
Set-TextCell $ws 25 5 '  -9.48%  ' $false

# Row 26
Set-TextCell $ws 26 4 '161.83' $true
Set-TextCell $ws 26 5 '  +6.02%  ' $false

# Row 27
Set-TextCell $ws 27 4 '20.74' $true
Set-TextCell $ws 27 5 '  +0.89%  ' $false

# Row 28
Set-TextCell $ws 28 4 '2.020.01' $false
Set-TextCell $ws 28 5 '  +3.68%  ' $false

# Row 29
Set-TextCell $ws 29 4 '2.282' $true
Set-TextCell $ws 29 5 '  -3.90%  ' $false

# Row 30
Set-TextCell $ws 30 4 '129.00' $true
Set-TextCell $ws 30 5 '  +1.35%  ' $false

# Row 31
Set-TextCell $ws 31 5 '  -2.14%  ' $false

# Row 32
Set-TextCell $ws 32 4 '6.031' $true
Set-TextCell $ws 32 5 '  +4.75%  ' $false

# Row 33
Set-TextCell $ws 33 4 '0.09159' $true
Set-TextCell $ws 33 5 '  +0.05%  ' $false

# Row 34
Set-TextCell $ws 34 4 '3.527' $true
Set-TextCell $ws 34 5 '  -3.45%  ' $false

# Row 35
Set-TextCell $ws 35 4 '12.95' $true
Set-TextCell $ws 35 5 '  +1.38%  ' $false

# Row 36
Set-TextCell $ws 36 4 '0.02361' $true
Set-TextCell $ws 36 5 '  +1.84%  ' $false

# Row 37
Set-TextCell $ws 37 4 '0.2175' $true
Set-TextCell $ws 37 5 '  +0.99%  ' $false

# Row 38
Set-TextCell $ws 38 4 '5.188' $true
Set-TextCell $ws 38 5 '  +1.21%  ' $false

# Row 39
Set-TextCell $ws 39 4 '0.6587' $true
Set-TextCell $ws 39 5 '  +1.36%  ' $false

# Row 40
Set-TextCell $ws 40 4 '0.06208' $true
Set-TextCell $ws 40 5 '  +1.44%  ' $false

# Row 41
Set-TextCell $ws 41 4 '1.190' $true
Set-TextCell $ws 41 5 '  -0.58%  ' $false

# Row 42
Set-TextCell $ws 42 4 '8.103' $true
Set-TextCell $ws 42 5 '  +1.52%  ' $false

# Row 43
Set-TextCell $ws 43 4 '1.427' $true
Set-TextCell $ws 43 5 '  +0.33%  ' $false

# Row 44
Set-TextCell $ws 44 4 '0.9985' $true
Set-TextCell $ws 44 5 '  +0.02%  ' $false

# Row 45
Set-TextCell $ws 45 4 '13.89' $true
Set-TextCell $ws 45 5 '  +0.77%  ' $false

# Row 46
Set-TextCell $ws 46 4 '0.6118' $true
Set-TextCell $ws 46 5 '  +2.86%  ' $false

# Row 47
Set-TextCell $ws 47 4 '3.737' $true
Set-TextCell $ws 47 5 '  -0.07%  ' $false

# Row 48
Set-TextCell $ws 48 4 '125.45' $true
Set-TextCell $ws 48 5 '  -0.68%  ' $false

# Row 49
Set-TextCell $ws 49 4 '2.018' $true
Set-TextCell $ws 49 5 '  +2.19%  ' $false

# Row 50 and 51 swap coins entirely (name, link, price, volume): Cronos <-> EOS
Set-TextCell $ws 50 2 'EOS' $false
Set-TextCell $ws 50 3 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos' $false
Set-TextCell $ws 50 4 '1.154' $true
Set-TextCell $ws 50 5 '  +2.47%  ' $false

Set-TextCell $ws 51 2 'Cronos' $false
Set-TextCell $ws 51 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' $false
Set-TextCell $ws 51 4 '0.06996' $true
Set-TextCell $ws 51 5 '  +1.35%  ' $false

